$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# "Queries" sheet: fix up the stray " =, " text in the SQL query string
# ---------------------------------------------------------------------------
$wsQueries = $wb.Worksheets.Item("Queries")
$wsQueries.Range("A2").Value = "SELECT [TemplateName] as 'Template Name'
      ,[FileName] as 'FileName'
      ,[LastChangedBy] as 'Last Changed By'
  ,Substring(LastChangedOn,11,18) as 'Last Changed On'
  FROM [Product_OCM].[dbo].[FaxTemplate];"
$wsQueries.Rows.Item(2).RowHeight = 75

# ---------------------------------------------------------------------------
# "Create" sheet: row 2 data re-typed (quote-prefixed so it stays plain text)
# ---------------------------------------------------------------------------
$wsCreate = $wb.Worksheets.Item("Create")
$wsCreate.Range("A2").Value = "'Mode"
$wsCreate.Range("C2").Value = "'Colors Group"
$wsCreate.Range("D2").Value = "'Mode Custom Template.html"
$wsCreate.Range("D2").Select()

# ---------------------------------------------------------------------------
# "Edit" sheet: row 2 data re-typed + page setup (paper size / orientation)
# ---------------------------------------------------------------------------
$wsEdit = $wb.Worksheets.Item("Edit")
$wsEdit.Range("A2").Value = "'Mode"
$wsEdit.Range("C2").Value = "'Colors Group"
$wsEdit.Range("D2").Value = "'Mode Custom Template.html"
$wsEdit.Range("E2").Value = "'Sachin Score"
$wsEdit.PageSetup.PaperSize = 9
$wsEdit.PageSetup.Orientation = 1

# ---------------------------------------------------------------------------
# "Delete" sheet: row 2 data re-typed, becomes the active sheet/tab
# ---------------------------------------------------------------------------
$wsDelete = $wb.Worksheets.Item("Delete")
$wsDelete.Range("A2").Value = "'Mode"
$wsDelete.Range("C2").Value = "'Colors Group"
$wsDelete.Range("D2").Value = "'Mode Custom Template.html"
$wsDelete.Range("E2").Value = "'Sachin Score"

# ---------------------------------------------------------------------------
# Make "Delete" the active sheet/tab (moves activeTab from Create to Delete)
# ---------------------------------------------------------------------------
$wsDelete.Activate()
